$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Create Sheet2 positioned after Sheet1
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Sheet2"

# Populate Sheet2 row 1 with new data (Laser Printers line item)
$ws2.Range("A1").Value = 1
$ws2.Range("B1").Value = "Laser Printers"
$ws2.Range("C1").Value = "SUPPLIES.OFFICE"
$ws2.Range("D1").Value = 500
$ws2.Range("E1").Value = "Each"
$ws2.Range("F1").Value = 100
$ws2.Range("G1").Value = "USD"
$ws2.Range("H1").Value = $false
$ws2.Range("I1").Value = "null"
$ws2.Range("J1").Value = $false
$ws2.Range("K1").Value = $false
$ws2.Range("L1").Value = "null"
$ws2.Range("M1").Value = "null"
$ws2.Range("N1").Value = "null"
$ws2.Range("O1").Value = "null"
$ws2.Range("P1").Value = "null"

# Sheet1: select the entire first row (matches sqref A1:XFD1) and drop tabSelected
$null = $ws1.Rows.Item(1).Select()

# Sheet2: set selection to H5 and make it the active/visible tab
$null = $ws2.Range("H5").Select()
$ws2.Activate()
